$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date values (stored as raw Excel serial numbers, cells already carry the
# date number format via their style) and Volumen/Precio promedio
# ponderado/Precio $/Kg updates, per the weekly re-shuffle of rows.

$updates = @{
    "D2"  = 44321
    "J2"  = 100
    "D3"  = 44321
    "J3"  = 50
    "D4"  = 44328
    "D5"  = 44328
    "D6"  = 44308
    "D7"  = 44308
    "D8"  = 44293
    "J8"  = 100
    "M8"  = 650
    "P8"  = 108
    "D9"  = 44293
    "D10" = 44188
    "D11" = 44188
    "D12" = 44525
    "J12" = 200
    "D13" = 44525
    "J13" = 100
    "D14" = 44230
    "D15" = 44230
    "D16" = 44358
    "D17" = 44358
    "D18" = 44335
    "J18" = 150
    "M18" = 633
    "P18" = 106
    "D19" = 44335
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
